$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035011844290189
$ws.Range("D2").Value = 1.043203085106238
$ws.Range("E2").Value = 1.034121493104194
$ws.Range("F2").Value = 1.053983195878501
$ws.Range("I2").Value = 1.039900491938862
$ws.Range("J2").Value = 1.040127949227614
$ws.Range("K2").Value = 1.045977613606845
$ws.Range("L2").Value = 1.036921876914136
$ws.Range("M2").Value = 1.056727660365422
$ws.Range("N2").Value = 1.041605050219861

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.035833716047721
$ws.Range("D3").Value = 1.043841260080336
$ws.Range("E3").Value = 1.034816962746946
$ws.Range("F3").Value = 1.0547880469562
$ws.Range("I3").Value = 1.040090341277645
$ws.Range("J3").Value = 1.040593836621374
$ws.Range("K3").Value = 1.04642722799138
$ws.Range("L3").Value = 1.037426771730944
$ws.Range("M3").Value = 1.057345685614662
$ws.Range("N3").Value = 1.042071599227158

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03636623632295
$ws.Range("D4").Value = 1.04425484182566
$ws.Range("E4").Value = 1.035267970143089
$ws.Range("F4").Value = 1.055309874068345
$ws.Range("I4").Value = 1.040212390554288
$ws.Range("J4").Value = 1.040895339892352
$ws.Range("K4").Value = 1.046718108567979
$ws.Range("L4").Value = 1.03775378748014
$ws.Range("M4").Value = 1.057745967264558
$ws.Range("N4").Value = 1.042373530667363

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036590277093489
$ws.Range("D5").Value = 1.044428862922902
$ws.Range("E5").Value = 1.035457809250725
$ws.Range("F5").Value = 1.055529495814483
$ws.Range("I5").Value = 1.040263508764662
$ws.Range("J5").Value = 1.041022100831209
$ws.Range("K5").Value = 1.046840381435254
$ws.Range("L5").Value = 1.037891338888606
$ws.Range("M5").Value = 1.057914334445614
$ws.Range("N5").Value = 1.042500471621293

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.036627904398419
$ws.Range("D6").Value = 1.044458090651095
$ws.Range("E6").Value = 1.035489697832655
$ws.Range("F6").Value = 1.055566385615659
$ws.Range("I6").Value = 1.04027208049693
$ws.Range("J6").Value = 1.041043385037756
$ws.Range("K6").Value = 1.046860910766699
$ws.Range("L6").Value = 1.037914438672299
$ws.Range("M6").Value = 1.057942609192121
$ws.Range("N6").Value = 1.042521786053855

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036369229303861
$ws.Range("D7").Value = 1.04425716651178
$ws.Range("E7").Value = 1.035270505859181
$ws.Range("F7").Value = 1.055312807702312
$ws.Range("I7").Value = 1.040213074350959
$ws.Range("J7").Value = 1.040897033644659
$ws.Range("K7").Value = 1.046719742438086
$ws.Range("L7").Value = 1.037755625159873
$ws.Range("M7").Value = 1.057748216648398
$ws.Range("N7").Value = 1.042375226824993

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035289450743289
$ws.Range("D8").Value = 1.043418625725016
$ws.Range("E8").Value = 1.034356323779609
$ws.Range("F8").Value = 1.054254983594799
$ws.Range("I8").Value = 1.039964816761568
$ws.Range("J8").Value = 1.040285388359162
$ws.Range("K8").Value = 1.046129572428346
$ws.Range("L8").Value = 1.037092442387944
$ws.Range("M8").Value = 1.056936445595698
$ws.Range("N8").Value = 1.041762712933034

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033392291356945
$ws.Range("D9").Value = 1.04194599224754
$ws.Range("E9").Value = 1.032753099788736
$ws.Range("F9").Value = 1.052398978025445
$ws.Range("I9").Value = 1.039521292718492
$ws.Range("J9").Value = 1.039207978948989
$ws.Range("K9").Value = 1.045089294659253
$ws.Range("L9").Value = 1.035926305879328
$ws.Range("M9").Value = 1.055508972126747
$ws.Range("N9").Value = 1.040683773477917

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032131353766648
$ws.Range("D10").Value = 1.04096769734228
$ws.Range("E10").Value = 1.031689561866119
$ws.Range("F10").Value = 1.051167153957617
$ws.Range("I10").Value = 1.039221579250554
$ws.Range("J10").Value = 1.038490045845459
$ws.Range("K10").Value = 1.0443956403104
$ws.Range("L10").Value = 1.035150632152065
$ws.Range("M10").Value = 1.054559423169078
$ws.Range("N10").Value = 1.039964820827033

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031586285694702
$ws.Range("D11").Value = 1.040544929163011
$ws.Range("E11").Value = 1.031230314120861
$ws.Range("F11").Value = 1.050635094344496
$ws.Range("I11").Value = 1.039090853534219
$ws.Range("J11").Value = 1.038179270226948
$ws.Range("K11").Value = 1.044095265282069
$ws.Range("L11").Value = 1.034815190228491
$ws.Range("M11").Value = 1.05414877780434
$ws.Range("N11").Value = 1.039653603871497

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031383964083542
$ws.Range("D12").Value = 1.040388022342949
$ws.Range("E12").Value = 1.031059922001665
$ws.Range("F12").Value = 1.050437665413808
$ws.Range("I12").Value = 1.039042154546817
$ws.Range("J12").Value = 1.038063849821582
$ws.Range("K12").Value = 1.043983691287947
$ws.Range("L12").Value = 1.034690658223719
$ws.Range("M12").Value = 1.053996325055291
$ws.Range("N12").Value = 1.039538019555916

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031427356349412
$ws.Range("D13").Value = 1.040421673560811
$ws.Range("E13").Value = 1.031096462927211
$ws.Range("F13").Value = 1.050480005428076
$ws.Range("I13").Value = 1.039052607042
$ws.Range("J13").Value = 1.038088607167385
$ws.Range("K13").Value = 1.04400762432654
$ws.Range("L13").Value = 1.034717367748464
$ws.Range("M13").Value = 1.054029023074188
$ws.Range("N13").Value = 1.039562812059991

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031569558841743
$ws.Range("D14").Value = 1.040531956571972
$ws.Range("E14").Value = 1.03121622550787
$ws.Range("F14").Value = 1.050618770680561
$ws.Range("I14").Value = 1.039086830945645
$ws.Range("J14").Value = 1.038169729216294
$ws.Range("K14").Value = 1.044086042560107
$ws.Range("L14").Value = 1.034804895017998
$ws.Range("M14").Value = 1.054136174386233
$ws.Range("N14").Value = 1.039644049311514

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031657193219898
$ws.Range("D15").Value = 1.040599922607236
$ws.Range("E15").Value = 1.031290040811635
$ws.Range("F15").Value = 1.050704295323355
$ws.Range("I15").Value = 1.039107898676905
$ws.Range("J15").Value = 1.038219713281865
$ws.Range("K15").Value = 1.044134358486617
$ws.Range("L15").Value = 1.034858832255248
$ws.Range("M15").Value = 1.054202204396542
$ws.Range("N15").Value = 1.039694104360191

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032167547806242
$ws.Range("D16").Value = 1.040995772934236
$ws.Range("E16").Value = 1.031720067576331
$ws.Range("F16").Value = 1.051202493170377
$ws.Range("I16").Value = 1.039230235181565
$ws.Range("J16").Value = 1.038510673098774
$ws.Range("K16").Value = 1.044415574948492
$ws.Range("L16").Value = 1.03517290351124
$ws.Range("M16").Value = 1.054586687362872
$ws.Range("N16").Value = 1.039985477373413

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032487928990785
$ws.Range("D17").Value = 1.041244305618708
$ws.Range("E17").Value = 1.031990153722888
$ws.Range("F17").Value = 1.051515356656364
$ws.Range("I17").Value = 1.039306720420974
$ws.Range("J17").Value = 1.038693210666092
$ws.Range("K17").Value = 1.044591970734285
$ws.Range("L17").Value = 1.035370028252487
$ws.Range("M17").Value = 1.054828002653222
$ws.Range("N17").Value = 1.040168274165014

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03267489116276
$ws.Range("D18").Value = 1.041389351445764
$ws.Range("E18").Value = 1.032147812923175
$ws.Range("F18").Value = 1.051697972656831
$ws.Range("I18").Value = 1.039351241466298
$ws.Range("J18").Value = 1.038799690734585
$ws.Range("K18").Value = 1.044694857570593
$ws.Range("L18").Value = 1.035485049161332
$ws.Range("M18").Value = 1.054968807409228
$ws.Range("N18").Value = 1.040274905447418

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032738655515528
$ws.Range("D19").Value = 1.041438821991282
$ws.Range("E19").Value = 1.032201591358165
$ws.Range("F19").Value = 1.051760261643708
$ws.Range("I19").Value = 1.039366406448967
$ws.Range("J19").Value = 1.038835999170644
$ws.Range("K19").Value = 1.044729938968431
$ws.Range("L19").Value = 1.035524275293253
$ws.Range("M19").Value = 1.055016826556838
$ws.Range("N19").Value = 1.040311265445621

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032453545877587
$ws.Range("D20").Value = 1.041217632030356
$ws.Range("E20").Value = 1.031961163348436
$ws.Range("F20").Value = 1.051481776107906
$ws.Range("I20").Value = 1.039298523737505
$ws.Range("J20").Value = 1.038673625172626
$ws.Range("K20").Value = 1.044573045316298
$ws.Range("L20").Value = 1.035348874342344
$ws.Range("M20").Value = 1.054802106660209
$ws.Range("N20").Value = 1.040148660857901

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03152767984707
$ws.Range("D21").Value = 1.040499477422255
$ws.Range("E21").Value = 1.031180953071021
$ws.Range("F21").Value = 1.050577902187322
$ws.Range("I21").Value = 1.039076756762746
$ws.Range("J21").Value = 1.038145840356168
$ws.Range("K21").Value = 1.044062950370536
$ws.Range("L21").Value = 1.03477911858326
$ws.Range("M21").Value = 1.054104618794422
$ws.Range("N21").Value = 1.039620126526466

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030946367447091
$ws.Range("D22").Value = 1.040048687238867
$ws.Range("E22").Value = 1.030691521170493
$ws.Range("F22").Value = 1.05001076880516
$ws.Range("I22").Value = 1.038936503981104
$ws.Range("J22").Value = 1.037814091422051
$ws.Range("K22").Value = 1.043742226398894
$ws.Range("L22").Value = 1.03442127346107
$ws.Range("M22").Value = 1.053666539723375
$ws.Range("N22").Value = 1.03928790647081

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031254454174401
$ws.Range("D23").Value = 1.040287588662741
$ws.Range("E23").Value = 1.030950871723801
$ws.Range("F23").Value = 1.050311305514454
$ws.Range("I23").Value = 1.039010932019129
$ws.Range("J23").Value = 1.037989948795667
$ws.Range("K23").Value = 1.043912248447947
$ws.Range("L23").Value = 1.034610937142115
$ws.Range("M23").Value = 1.053898729609898
$ws.Range("N23").Value = 1.039464013582068

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032469081850903
$ws.Range("D24").Value = 1.041229684427016
$ws.Range("E24").Value = 1.03197426247247
$ws.Range("F24").Value = 1.051496949317011
$ws.Range("I24").Value = 1.039302227748762
$ws.Range("J24").Value = 1.038682474986462
$ws.Range("K24").Value = 1.044581596903921
$ws.Range("L24").Value = 1.035358432756107
$ws.Range("M24").Value = 1.054813807791654
$ws.Range("N24").Value = 1.040157523239487

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03388208422228
$ws.Range("D25").Value = 1.04232610142354
$ws.Range("E25").Value = 1.033166649616652
$ws.Range("F25").Value = 1.052877837308705
$ws.Range("I25").Value = 1.039636668088423
$ws.Range("J25").Value = 1.039486461004465
$ws.Range("K25").Value = 1.045358261285654
$ws.Range("L25").Value = 1.036227477834607
$ws.Range("M25").Value = 1.055877645603506
$ws.Range("N25").Value = 1.040962651009855

